$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 28 ("Таблица - DE_TAB_ListInventoryCommission"
# title row). This shifts every row from 28 onward down by one (data, titles, blank
# separators, merged-cell title bars, and the trailing blank row all move together), and
# a fresh blank row appears at the very bottom (old last row 1003 -> 1004).
$ws.Rows(28).Insert()

# The freshly inserted row 28 is blank; give it the same formatting as a standard
# attribute row of a table (the "de_employeeId" row, now shifted to row 34, has exactly
# the style pattern we need: s=5,5,5,5,8,5,7,5).
$ws.Range("A34:H34").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)

# Populate the new row: a 4th attribute added to the "DE_CTL_Employees" table,
# de_employeePositionId / Должность сотрудника / FK / INTEGER / NOT NULL / >0
$ws.Range("A28").Value = "de_employeePositionId"
$ws.Range("B28").Value = "Должность сотрудника"
$ws.Range("C28").Value = "FK"
$ws.Range("D28").Value = "INTEGER"
$ws.Range("F28").Value = "NOT NULL"
$ws.Range("H28").Value = ">0"
